# Insert a new weekly price record at row 62 ("Región de O'Higgins", 2022-01-11),
# pushing the existing rows 62-89 down to 63-90.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(62).Insert()

$ws.Cells.Item(62, 1).Value2  = 11
$ws.Cells.Item(62, 2).Value2  = 'Vega Monumental Concepción'
$ws.Cells.Item(62, 3).Value2  = 'Bíobío'
$ws.Cells.Item(62, 4).Value2  = 44572
$ws.Cells.Item(62, 5).Value2  = 8
$ws.Cells.Item(62, 6).Value2  = 100112032
$ws.Cells.Item(62, 7).Value2  = 'Zapallo italiano'
$ws.Cells.Item(62, 8).Value2  = 'Sin especificar'
$ws.Cells.Item(62, 9).Value2  = 'Primera'
$ws.Cells.Item(62, 10).Value2 = 260
$ws.Cells.Item(62, 11).Value2 = 14000
$ws.Cells.Item(62, 12).Value2 = 15000
$ws.Cells.Item(62, 13).Value2 = 14538
$ws.Cells.Item(62, 14).Value2 = '$/caja 60 unidades'
$ws.Cells.Item(62, 15).Value2 = 'Región de O''Higgins'
$ws.Cells.Item(62, 16).Value2 = 242
$ws.Cells.Item(62, 17).Value2 = 60
$ws.Cells.Item(62, 18).Value2 = 'Hortaliza'
